# Week 2 Materials updated
#
# Two visible text edits:
#   1. Slide 23 - "Text Placeholder 2": merge the two runs
#      "IQR(): inner quartile " + "range (Q3 - Q1)" into a single run
#      "IQR(): inner quartile range (Q3 - Q1)".
#   2. Slide 26 - "Text Placeholder 2", paragraph 7: change
#      "Bounded between 0 and 1." to "Bounded between -1 and 1."

$p = $ppt.ActivePresentation

# --- Slide 23: IQR() placeholder -------------------------------------------
$s23 = $p.Slides.Item(23)
$sh23 = $s23.Shapes.Item(2)              # "Text Placeholder 2"
$tr23 = $sh23.TextFrame.TextRange
$iqrPara = $tr23.Paragraphs(1)

$run1 = $iqrPara.Runs(1)
$run2 = $iqrPara.Runs(2)
$run1.Text = "IQR(): inner quartile range (Q3 " + [char]0x2013 + " Q1)"
$run2.Text = ""

# --- Slide 26: correlation bounds placeholder -------------------------------
$s26 = $p.Slides.Item(26)
$sh26 = $s26.Shapes.Item(3)              # "Text Placeholder 2"
$tr26 = $sh26.TextFrame.TextRange
$boundedPara = $tr26.Paragraphs(7)

$b1 = $boundedPara.Runs(1)
$b1.Text = "Bounded "
$b2 = $b1.InsertAfter("between ")
$b3 = $b2.InsertAfter("-1 ")
$b4 = $b3.InsertAfter("and 1.")
